{"js": "const body = context.document.body;\n\n// The doc reads \"Version 1.\" -> bump the version number to \"Version 2.\"\n// (the trailing period sits after the \"_GoBack\" bookmark in the source,\n// so it is appended separately to preserve that layout).\n\n// 1) \" 1.\" -> \" 2\"  (drop the old period along with the old number)\nlet numberMatches = body.search(\" 1.\", { matchCase: true, matchWholeWord: false });\nnumberMatches.load(\"items\");\nawait context.sync();\n\nif (numberMatches.items.length === 0) {\n  throw new Error('Could not find \"Version 1.\" text to update.');\n}\nnumberMatches.items[0].insertText(\" 2\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) Re-append the final period at the very end of the body (after the bookmark).\nbody.insertText(\".\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The doc reads \"Version 1.\" -> bump the version number to \"Version 2.\"\n# (the trailing period sits after the \"_GoBack\" bookmark in the source,\n# so it is appended separately to preserve that layout).\n\n# 1) \" 1.\" -> \" 2\"  (drop the old period along with the old number)\n$r = $d.Content\n$find = $r.Find\n$find.ClearFormatting()\n$find.Text = \" 1.\"\n$found = $find.Execute()\nif (-not $found) {\n  throw \"Could not find 'Version 1.' text to update.\"\n}\n$r.Text = \" 2\"\n\n# 2) Re-append the final period at the very end of the document (after the bookmark).\n$endRange = $d.Range($d.Content.End, $d.Content.End)\n$endRange.InsertAfter(\".\")\n"}
